# "Generate Report for Archive" — refresh the localization-status report:
# the zh-cn/de-de handoff moved from "Ready for handoff" to "In Translation",
# so update the Overview rollup (zh-cn + de-de status columns) and each
# locale sheet's own Status column, then let the (now shorter) Status
# columns re-fit to their content, same as the reporting job does on
# every regeneration.

$wb = $excel.ActiveWorkbook

# Overview sheet: columns E ("zh-cn") and F ("de-de") hold that locale's status.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Columns.Item(5).AutoFit() | Out-Null
$overview.Columns.Item(6).AutoFit() | Out-Null
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

# Per-locale detail sheets: column C is "Status".
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Columns.Item(3).AutoFit() | Out-Null
$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Columns.Item(3).AutoFit() | Out-Null
$dede.Columns.Item(3).ColumnWidth = 12.5
